# Applies the "updated the chapters with the tabu search" revision:
#  1. Cover-page date "DESAMBER 2022" -> "JULY 2023"
#  2. Drop the stray "s" that turned "account" into "accounts" in the
#     problem-statement paragraph.
#  3. "Generate university lectures" -> "Generate university lecture"
#     (singular) in the objectives list.
#  4. Capitalize "university" -> "University" right after "Hadhramout".
#  5. Capitalize "college of engineering and petroleum" ->
#     "College of engineering and Petroleum".

$d = $word.ActiveDocument

# --- 1. DESAMBER 2022 -> JULY 2023 -------------------------------------
$r = $d.Content
$r.Find.Execute("DESAMBER 2022", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
if ($r.Find.Found) {
    $r.Text = "JULY 2023"
}

# --- 2. "...taken into accounts, such as..." -> "...account, such as..." -
$r = $d.Content
$r.Find.Execute("taken into accounts,", $true, $false, $false, $false, `
                 $false, $true, 1, $false, "", 0) | Out-Null
if ($r.Find.Found) {
    $r.Text = "taken into account,"
}

# --- 3. "Generate university lectures" -> "Generate university lecture" -
$r = $d.Content
$r.Find.Execute("Generate university lectures", $true, $false, $false, `
                 $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r.Find.Found) {
    $r.Text = "Generate university lecture"
}

# --- 4. " university" -> " University" (only the one after "Hadhramout") -
# Scope the search to start after the "Generate university lectures" list
# item so the earlier (lower-case, unrelated) occurrence is not touched.
$scope = $d.Range(2300, $d.Content.End)
$scope.Find.Execute(" university", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
if ($scope.Find.Found) {
    $scope.Text = " University"
}

# --- 5. "college of engineering and petroleum" -> "College of engineering
#         and Petroleum" ---------------------------------------------------
$scope2 = $d.Range(2300, $d.Content.End)
$scope2.Find.Execute("college of engineering and petroleum", $true, $false, `
                      $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($scope2.Find.Found) {
    $scope2.Text = "College of engineering and Petroleum"
}
